$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1237.01
$ws.Range("I112").Value = 399.0909
$ws.Range("J112").Value = 1340.573
$ws.Range("K112").Value = 1197.2727
$ws.Range("L112").Value = 4021.719
$ws.Range("M112").Value = -89.27269999999999
$ws.Range("N112").Value = -6237.719

$ws.Range("H132").Value = 9531109
$ws.Range("I132").Value = 16676156
$ws.Range("J132").Value = 4379.1113
$ws.Range("K132").Value = 50028468
$ws.Range("L132").Value = 13137.3339
$ws.Range("M132").Value = -50025938
$ws.Range("N132").Value = -18197.3339

$ws.Range("H137").Value = 7695088.5
$ws.Range("I137").Value = 11113456
$ws.Range("J137").Value = 3762.25
$ws.Range("K137").Value = 33340368
$ws.Range("L137").Value = 11286.75
$ws.Range("M137").Value = -33337818
$ws.Range("N137").Value = -16386.75

$ws.Range("H138").Value = 3177.9185
$ws.Range("I138").Value = 1560.7
$ws.Range("J138").Value = 5731.421
$ws.Range("K138").Value = 4682.1
$ws.Range("L138").Value = 17194.263
$ws.Range("M138").Value = 457.8999999999996
$ws.Range("N138").Value = -27474.263

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 83338340
$ws.Range("I2").Value = 125002500
$ws.Range("J2").Value = 10000
$ws.Range("K2").Value = 125002500
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = -125002387
$ws.Range("N2").Value = -10226

$ws.Range("H45").Value = 1694.6765
$ws.Range("I45").Value = 1088.5862
$ws.Range("K45").Value = 1088.5862
$ws.Range("M45").Value = -711.5862

$ws.Range("H74").Value = 1004.6667
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 1007
$ws.Range("K74").Value = 1000
$ws.Range("L74").Value = 1007
$ws.Range("M74").Value = -126
$ws.Range("N74").Value = -2755

$ws.Range("H77").Value = 1004.6667
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 1007
$ws.Range("K77").Value = 5000
$ws.Range("L77").Value = 5035
$ws.Range("M77").Value = -632
$ws.Range("N77").Value = -13771

$ws.Range("H88").Value = 1826.2858
$ws.Range("I88").Value = 1826.2858
$ws.Range("K88").Value = 1826.2858
$ws.Range("M88").Value = -1420.2858

$ws.Range("H91").Value = 1826.2858
$ws.Range("I91").Value = 1826.2858
$ws.Range("K91").Value = 1826.2858
$ws.Range("M91").Value = -422.2858000000001

$ws.Range("H97").Value = 693.93335
$ws.Range("I97").Value = 693.93335
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 693.93335
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -197.93335
$ws.Range("N97").ClearContents()

$ws.Range("H102").Value = 2379.2354
$ws.Range("I102").Value = 2371.6875
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 2371.6875
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -749.6875
$ws.Range("N102").Value = -5744

$ws.Range("H116").Value = 83338340
$ws.Range("I116").Value = 125002500
$ws.Range("J116").Value = 10000
$ws.Range("K116").Value = 125002500
$ws.Range("L116").Value = 10000
$ws.Range("M116").Value = -125000206
$ws.Range("N116").Value = -14588

$ws.Range("H122").Value = 3072.182
$ws.Range("I122").Value = 2203.4546
$ws.Range("J122").Value = 3940.9092
$ws.Range("K122").Value = 6610.3638
$ws.Range("L122").Value = 11822.7276
$ws.Range("M122").Value = -4160.3638
$ws.Range("N122").Value = -16722.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 83338340
$ws.Range("I3").Value = 125002500
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 125002500
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = -125002386
$ws.Range("N3").Value = -10228

$ws.Range("H94").Value = 602.6667
$ws.Range("I94").Value = 578.7
$ws.Range("K94").Value = 578.7
$ws.Range("M94").Value = -127.7

$ws.Range("H99").Value = 1816.3478
$ws.Range("I99").Value = 1320.375
$ws.Range("J99").Value = 2950
$ws.Range("K99").Value = 1320.375
$ws.Range("L99").Value = 2950
$ws.Range("M99").Value = 177.625
$ws.Range("N99").Value = -5946

$ws.Range("H134").Value = 3579.0667
$ws.Range("I134").Value = 3429.2173
$ws.Range("J134").Value = 4071.4285
$ws.Range("K134").Value = 10287.6519
$ws.Range("L134").Value = 12214.2855
$ws.Range("M134").Value = -7752.651899999999
$ws.Range("N134").Value = -17284.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1889161
$ws.Range("I31").Value = 2175522
$ws.Range("J31").Value = 7360.7144
$ws.Range("K31").Value = 2175522
$ws.Range("L31").Value = 7360.7144
$ws.Range("M31").Value = -2175227
$ws.Range("N31").Value = -7950.7144

$ws.Range("H34").Value = 1889161
$ws.Range("I34").Value = 2175522
$ws.Range("J34").Value = 7360.7144
$ws.Range("K34").Value = 2175522
$ws.Range("L34").Value = 7360.7144
$ws.Range("M34").Value = -2175320
$ws.Range("N34").Value = -7764.7144

$ws.Range("H122").Value = 1898.3125
$ws.Range("I122").Value = 1947.1875
$ws.Range("J122").Value = 1849.4375
$ws.Range("K122").Value = 5841.5625
$ws.Range("L122").Value = 5548.3125
$ws.Range("M122").Value = -3391.5625
$ws.Range("N122").Value = -10448.3125

$ws.Range("H134").Value = 2257.6428
$ws.Range("I134").Value = 1087.4286
$ws.Range("J134").Value = 3427.8572
$ws.Range("K134").Value = 3262.2858
$ws.Range("L134").Value = 10283.5716
$ws.Range("M134").Value = -727.2857999999997
$ws.Range("N134").Value = -15353.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1532.5
$ws.Range("I97").Value = 965.7143
$ws.Range("J97").Value = 5500
$ws.Range("K97").Value = 965.7143
$ws.Range("L97").Value = 5500
$ws.Range("M97").Value = -469.7143
$ws.Range("N97").Value = -6492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 695.3913
$ws.Range("I55").Value = 171.83333
$ws.Range("J55").Value = 880.17645
$ws.Range("K55").Value = 171.83333
$ws.Range("L55").Value = 880.17645
$ws.Range("M55").Value = 1.166670000000011
$ws.Range("N55").Value = -1226.17645

$ws.Range("H122").Value = 2745.4546
$ws.Range("I122").Value = 2460.1667
$ws.Range("J122").Value = 3356.7856
$ws.Range("K122").Value = 7380.500100000001
$ws.Range("L122").Value = 10070.3568
$ws.Range("M122").Value = -4930.500100000001
$ws.Range("N122").Value = -14970.3568

$ws.Range("H132").Value = 3291.8333
$ws.Range("I132").Value = 1615.8462
$ws.Range("J132").Value = 5272.5454
$ws.Range("K132").Value = 4847.5386
$ws.Range("L132").Value = 15817.6362
$ws.Range("M132").Value = -2317.5386
$ws.Range("N132").Value = -20877.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 19749.666
$ws.Range("J63").Value = 19749.666
$ws.Range("L63").Value = 19749.666
$ws.Range("N63").Value = -20997.666

$ws.Range("H66").Value = 19749.666
$ws.Range("J66").Value = 19749.666
$ws.Range("L66").Value = 59248.99800000001
$ws.Range("N66").Value = -65488.99800000001

$ws.Range("H69").Value = 26799.875
$ws.Range("J69").Value = 26799.875
$ws.Range("L69").Value = 26799.875
$ws.Range("N69").Value = -28297.875

$ws.Range("H72").Value = 26799.875
$ws.Range("J72").Value = 26799.875
$ws.Range("L72").Value = 80399.625
$ws.Range("N72").Value = -87887.625

$ws.Range("H75").Value = 31606
$ws.Range("J75").Value = 31606
$ws.Range("L75").Value = 31606
$ws.Range("N75").Value = -33478

$ws.Range("H76").Value = 30900
$ws.Range("I76").Value = 20000
$ws.Range("J76").Value = 32111.111
$ws.Range("K76").Value = 20000
$ws.Range("L76").Value = 32111.111
$ws.Range("M76").Value = -19685
$ws.Range("N76").Value = -32741.111

$ws.Range("H78").Value = 31606
$ws.Range("J78").Value = 31606
$ws.Range("L78").Value = 94818
$ws.Range("N78").Value = -104178

$ws.Range("H79").Value = 30900
$ws.Range("I79").Value = 20000
$ws.Range("J79").Value = 32111.111
$ws.Range("K79").Value = 20000
$ws.Range("L79").Value = 32111.111
$ws.Range("M79").Value = -18908
$ws.Range("N79").Value = -34295.111

$ws.Range("H80").Value = 38333.332
$ws.Range("J80").Value = 38333.332
$ws.Range("L80").Value = 38333.332
$ws.Range("N80").Value = -40329.332

$ws.Range("H82").Value = 31636.363
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 31636.363
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 31636.363
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -32402.363

$ws.Range("H83").Value = 38333.332
$ws.Range("J83").Value = 38333.332
$ws.Range("L83").Value = 114999.996
$ws.Range("N83").Value = -124983.996

$ws.Range("H85").Value = 31636.363
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 31636.363
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 31636.363
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -34288.363

$ws.Range("H126").Value = 4764081
$ws.Range("I126").Value = 1664
$ws.Range("J126").Value = 14288915
$ws.Range("K126").Value = 4992
$ws.Range("L126").Value = 42866745
$ws.Range("M126").Value = -2522
$ws.Range("N126").Value = -42871685

$ws.Range("H138").Value = 29429
$ws.Range("J138").Value = 29429
$ws.Range("L138").Value = 29429
$ws.Range("N138").Value = -39709
